$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2-117)
# from 45179 (2023-09-10) to 45180 (2023-09-11).
for ($r = 2; $r -le 117; $r++) {
    $ws.Cells.Item($r, 3).Value = 45180
}

# Row 107 (Beteckning "A 7181-2023") gains new HYPERLINK formulas in
# columns U through Y.
$ws.Range("U107").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MULLSJO/knärot/A 7181-2023.png")'
$ws.Range("V107").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MULLSJO/klagomål/A 7181-2023.docx")'
$ws.Range("W107").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MULLSJO/klagomålsmail/A 7181-2023.docx")'
$ws.Range("X107").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MULLSJO/tillsyn/A 7181-2023.docx")'
$ws.Range("Y107").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MULLSJO/tillsynsmail/A 7181-2023.docx")'
